$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 58
$ws.Range("B1").Value = 43
$ws.Range("C1").Value = 40

$ws.Range("A2").Value = 142
$ws.Range("B2").Value = 115
$ws.Range("C2").Value = 115

$ws.Range("A3").Value = 226
$ws.Range("B3").Value = 187
$ws.Range("C3").Value = 190
